$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Update TestCases value for row 2 from "26,28,29,30,31,36" to "34"
$ws.Range("B2").Value = "34"

# Move the active selection from D2 to C2
$ws.Range("C2").Select()
